$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value2 = 150
$ws.Range("J5").Value2 = 0
$ws.Range("L5").Value2 = 0
$ws.Range("N5").ClearContents()
$ws.Range("H17").Value2 = 3994
$ws.Range("J17").Value2 = 3994
$ws.Range("L17").Value2 = 11982
$ws.Range("N17").Value2 = -12318
$ws.Range("H43").Value2 = 4500
$ws.Range("I43").Value2 = 7000
$ws.Range("K43").Value2 = 7000
$ws.Range("M43").Value2 = -6931
$ws.Range("H70").Value2 = 10000000
$ws.Range("J70").Value2 = 10000000
$ws.Range("L70").Value2 = 30000000
$ws.Range("N70").Value2 = -30000540
$ws.Range("H73").Value2 = 10000000
$ws.Range("J73").Value2 = 10000000
$ws.Range("L73").Value2 = 30000000
$ws.Range("N73").Value2 = -30001872
$ws.Range("H92").Value2 = 348.9
$ws.Range("I92").Value2 = 343.44446
$ws.Range("J92").Value2 = 398
$ws.Range("K92").Value2 = 343.44446
$ws.Range("L92").Value2 = 398
$ws.Range("M92").Value2 = 904.5555400000001
$ws.Range("N92").Value2 = -2894
$ws.Range("H107").Value2 = 277.875
$ws.Range("I107").Value2 = 289.14285
$ws.Range("K107").Value2 = 289.14285
$ws.Range("M107").Value2 = 1630.85715
$ws.Range("H112").Value2 = 1240.25
$ws.Range("J112").Value2 = 1730.9474
$ws.Range("L112").Value2 = 5192.8422
$ws.Range("N112").Value2 = -7408.8422
$ws.Range("H113").Value2 = 5001503
$ws.Range("J113").Value2 = 3006
$ws.Range("L113").Value2 = 3006
$ws.Range("N113").Value2 = -9514
$ws.Range("H116").Value2 = 6997.1665
$ws.Range("J116").Value2 = 6994.75
$ws.Range("L116").Value2 = 6994.75
$ws.Range("N116").Value2 = -13878.75
$ws.Range("H132").Value2 = 2697.9
$ws.Range("I132").Value2 = 1920.4615
$ws.Range("K132").Value2 = 5761.3845
$ws.Range("M132").Value2 = -3231.3845
$ws.Range("H137").Value2 = 1780.6666
$ws.Range("I137").Value2 = 1628
$ws.Range("K137").Value2 = 4884
$ws.Range("M137").Value2 = -2334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 851.2778
$ws.Range("I2").Value2 = 1044.5385
$ws.Range("K2").Value2 = 1044.5385
$ws.Range("M2").Value2 = -931.5385000000001
$ws.Range("H32").Value2 = 12679.272
$ws.Range("I32").Value2 = 12679.272
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 12679.272
$ws.Range("L32").Value2 = 0
$ws.Range("M32").Value2 = -12392.272
$ws.Range("N32").ClearContents()
$ws.Range("H45").Value2 = 1999
$ws.Range("I45").Value2 = 1999
$ws.Range("K45").Value2 = 1999
$ws.Range("M45").Value2 = -1622
$ws.Range("H61").Value2 = 2454.423
$ws.Range("I61").Value2 = 1916.4286
$ws.Range("J61").Value2 = 3082.0833
$ws.Range("K61").Value2 = 1916.4286
$ws.Range("L61").Value2 = 3082.0833
$ws.Range("M61").Value2 = -1704.4286
$ws.Range("N61").Value2 = -3506.0833
$ws.Range("H102").Value2 = 1485.2858
$ws.Range("I102").Value2 = 1485.2858
$ws.Range("K102").Value2 = 1485.2858
$ws.Range("M102").Value2 = 136.7141999999999
$ws.Range("H116").Value2 = 851.2778
$ws.Range("I116").Value2 = 1044.5385
$ws.Range("K116").Value2 = 1044.5385
$ws.Range("M116").Value2 = 1249.4615
$ws.Range("H122").Value2 = 3508.7727
$ws.Range("I122").Value2 = 3620.8
$ws.Range("J122").Value2 = 2388.5
$ws.Range("K122").Value2 = 10862.4
$ws.Range("L122").Value2 = 7165.5
$ws.Range("M122").Value2 = -8412.400000000001
$ws.Range("N122").Value2 = -12065.5
$ws.Range("H132").Value2 = 2368.6296
$ws.Range("I132").Value2 = 1558.7222
$ws.Range("J132").Value2 = 3988.4443
$ws.Range("K132").Value2 = 4676.1666
$ws.Range("L132").Value2 = 11965.3329
$ws.Range("M132").Value2 = -2146.1666
$ws.Range("N132").Value2 = -17025.3329
$ws.Range("H136").Value2 = 2454.423
$ws.Range("I136").Value2 = 1916.4286
$ws.Range("J136").Value2 = 3082.0833
$ws.Range("K136").Value2 = 5749.2858
$ws.Range("L136").Value2 = 9246.249899999999
$ws.Range("M136").Value2 = -3199.2858
$ws.Range("N136").Value2 = -14346.2499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 851.2778
$ws.Range("I3").Value2 = 1044.5385
$ws.Range("K3").Value2 = 1044.5385
$ws.Range("M3").Value2 = -930.5385000000001
$ws.Range("H22").Value2 = 680.4
$ws.Range("I22").Value2 = 675.875
$ws.Range("J22").Value2 = 698.5
$ws.Range("K22").Value2 = 675.875
$ws.Range("L22").Value2 = 698.5
$ws.Range("M22").Value2 = -502.875
$ws.Range("N22").Value2 = -1044.5
$ws.Range("H94").Value2 = 438.8889
$ws.Range("I94").Value2 = 438.8889
$ws.Range("J94").Value2 = 0
$ws.Range("K94").Value2 = 438.8889
$ws.Range("L94").Value2 = 0
$ws.Range("M94").Value2 = 12.11110000000002
$ws.Range("N94").ClearContents()
$ws.Range("H134").Value2 = 1575.7715
$ws.Range("I134").Value2 = 1395.2258
$ws.Range("J134").Value2 = 2975
$ws.Range("K134").Value2 = 4185.6774
$ws.Range("L134").Value2 = 8925
$ws.Range("M134").Value2 = -1650.6774
$ws.Range("N134").Value2 = -13995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1380
$ws.Range("J16").Value2 = 997
$ws.Range("L16").Value2 = 997
$ws.Range("N16").Value2 = -1571
$ws.Range("H105").Value2 = 1628.6
$ws.Range("I105").Value2 = 1665.1111
$ws.Range("K105").Value2 = 1665.1111
$ws.Range("M105").Value2 = 81.88889999999992
$ws.Range("H107").Value2 = 1457.1875
$ws.Range("I107").Value2 = 691.5
$ws.Range("K107").Value2 = 691.5
$ws.Range("M107").Value2 = 1228.5
$ws.Range("H113").Value2 = 1380
$ws.Range("J113").Value2 = 997
$ws.Range("L113").Value2 = 997
$ws.Range("N113").Value2 = -5337
$ws.Range("H132").Value2 = 2536.4583
$ws.Range("I132").Value2 = 1929.7333
$ws.Range("J132").Value2 = 3547.6667
$ws.Range("K132").Value2 = 5789.199900000001
$ws.Range("L132").Value2 = 10643.0001
$ws.Range("M132").Value2 = -3259.199900000001
$ws.Range("N132").Value2 = -15703.0001
$ws.Range("H134").Value2 = 2699.15
$ws.Range("I134").Value2 = 2734
$ws.Range("J134").Value2 = 2594.6
$ws.Range("K134").Value2 = 8202
$ws.Range("L134").Value2 = 7783.799999999999
$ws.Range("M134").Value2 = -5667
$ws.Range("N134").Value2 = -12853.8
$ws.Range("H141").Value2 = 60000
$ws.Range("J141").Value2 = 60000
$ws.Range("L141").Value2 = 60000
$ws.Range("N141").Value2 = -70360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 169531280
$ws.Range("I4").Value2 = 105038950
$ws.Range("M4").Value2 = -315116738
$ws.Range("H8").Value2 = 375
$ws.Range("I8").Value2 = 375
$ws.Range("M8").Value2 = -986
$ws.Range("H113").Value2 = 1892.3684
$ws.Range("J113").Value2 = 1903.1875
$ws.Range("L113").Value2 = 5709.5625
$ws.Range("N113").Value2 = -10049.5625
$ws.Range("H137").Value2 = 2986.875
$ws.Range("I137").Value2 = 2779
$ws.Range("K137").Value2 = 8337
$ws.Range("M137").Value2 = -3237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 1073.5
$ws.Range("I102").Value2 = 1058.75
$ws.Range("K102").Value2 = 1058.75
$ws.Range("M102").Value2 = 563.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 277.75
$ws.Range("I16").Value2 = 277.75
$ws.Range("K16").Value2 = 277.75
$ws.Range("M16").Value2 = -107.75
$ws.Range("H43").Value2 = 174832.67
$ws.Range("J43").Value2 = 174832.67
$ws.Range("L43").Value2 = 174832.67
$ws.Range("N43").Value2 = -175218.67
$ws.Range("H61").Value2 = 4999
$ws.Range("I61").Value2 = 0
$ws.Range("J61").Value2 = 4999
$ws.Range("K61").Value2 = 0
$ws.Range("L61").Value2 = 4999
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value2 = -5403
$ws.Range("H113").Value2 = 4999
$ws.Range("I113").Value2 = 0
$ws.Range("J113").Value2 = 4999
$ws.Range("K113").Value2 = 0
$ws.Range("L113").Value2 = 4999
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value2 = -9339
